$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the extraction-template values in E2 and E10
$ws.Range("E2").Value = "\Testdata\Templates\ManagePopulations\extraction-template-12.xlsx"
$ws.Range("E10").Value = "\Testdata\Templates\ManagePopulations\extraction-template-17.xlsx"

# Update the current selection on the sheet
$ws.Range("E2:E10").Select()
